$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.148.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "'2.814.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'361.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").Value = "'110.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.13%  "

$ws.Range("D7").Value = "'0.564"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.30%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").Value = "'40.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.31%  "

$ws.Range("D11").Value = "'0.0859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "'19.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").Value = "'7.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("D15").Value = "'3.259.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").Value = "'2.844.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.93%  "

$ws.Range("D17").Value = "'0.919"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("D18").Value = "'52.039.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("E19").Value = "  +1.46%  "

$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("D21").Value = "'13.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.84%  "

$ws.Range("D22").Value = "'0.0₃0990"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'272.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'69.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").Value = "'26.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.64%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "'10.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.86%  "

$ws.Range("E29").Value = "  -1.18%  "

$ws.Range("E30").Value = "  +0.62%  "

$ws.Range("D31").Value = "'0.0475"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.20%  "

$ws.Range("D32").Value = "'52.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.48%  "

$ws.Range("D33").Value = "'34.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.39%  "

$ws.Range("D34").Value = "'5.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("D35").Value = "'5.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.66%  "

$ws.Range("D36").Value = "'0.0844"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "

$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D38").Value = "'3.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.02%  "

$ws.Range("D39").Value = "'2.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.49%  "

$ws.Range("D40").Value = "'18.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.80%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").Value = "'2.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.07%  "

$ws.Range("D43").Value = "'125.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("D44").Value = "'2.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("D45").Value = "'22.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.07%  "

$ws.Range("D46").Value = "'2.068.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").Value = "'3.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.08%  "

$ws.Range("D48").Value = "'2.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.16%  "

$ws.Range("D49").Value = "'5.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.82%  "

$ws.Range("D50").Value = "'0.949"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.83%  "

$ws.Range("D51").Value = "'9.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.99%  "
